# docs: add MongoDB database creation subtasks to project plan
#
# Insert a new task row (row 3) on the project tracker sheet for the new
# "MongoDB Database Architecture" task, shifting the existing "Security
# Testing" / "Vulnerability Engine" / "Infrastructure" / "Observability" /
# "Production Security" / "CI/CD" rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3 ("Linear Jailbreaking
# Strategy"), pushing it (and everything below) down to row 4.
$ws.Rows.Item(3).Insert()

# Populate the newly-inserted row 3 with the MongoDB Database Architecture task.
$ws.Range("A3").Value = "Core Framework"
$ws.Range("B3").Value = "Low-Level"
$ws.Range("C3").Value = "MongoDB Database Architecture"
$ws.Range("D3").Value = "Design and implement a scalable document storage system for red-teaming runs, results, and profiles."
$ws.Range("E3").Value = "Database Engineer"
$ws.Range("F3").Value = "Critical"
$ws.Range("G3").Value = "Implemented"
$ws.Range("H3").Value = "Medium"
$ws.Range("I3").Value = "1. Design and document document schema (Collections for Runs, Results, and Profiles)"
$ws.Range("J3").Value = "2. Setup MongoDB environment (Atlas/Local) and configure connection strings"
$ws.Range("K3").Value = "3. Implement StorageHelper and MongoDBService classes for CRUD operations"
$ws.Range("L3").Value = "4. Develop data migration scripts and validation logic for schema consistency"
